$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.368870615959167
$ws.Range("B1").Value = 2.752486944198608
$ws.Range("C1").Value = 3.449557781219482
$ws.Range("D1").Value = 3.476047039031982
$ws.Range("E1").Value = 1.570347547531128
